$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename existing sheets
# ---------------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Invoices")
$wsSales.Name = "invoice_sales"

$wsSettlement = $wb.Worksheets.Item("Settlement")
$wsSettlement.Name = "settlement"

# ---------------------------------------------------------------------------
# 2. Add the new "invoice_costs" sheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCosts = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsCosts.Name = "invoice_costs"

# ---------------------------------------------------------------------------
# 3. Populate "invoice_costs" - mirror the layout used on "invoice_sales"
#    (Field name >>> / Value type >>> / Field expression >>> / Search values >>>)
# ---------------------------------------------------------------------------

# Pull formats first (keeps style indexes identical instead of minting new ones):
#   A1:E4 -> plain "Arial 10" body style (same as invoice_sales!A2)
#   B1:E1 -> bold/centered/bordered header style (same as invoice_sales!B1)
$wsSales.Range("A2").Copy() | Out-Null
$wsCosts.Range("A1:E4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsSales.Range("B1").Copy() | Out-Null
$wsCosts.Range("B1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 1 - field names
$wsCosts.Range("A1").Value = "Field name >>>"
$wsCosts.Range("B1").Value = "Invoice Number"
$wsCosts.Range("C1").Value = "Invoice Date"
$wsCosts.Range("D1").Value = "Product"
$wsCosts.Range("E1").Value = "Total"

# Row 2 - value types
$wsCosts.Range("A2").Value = "Value type >>>"
$wsCosts.Range("B2").Value = "String"
$wsCosts.Range("C2").Value = "String"
$wsCosts.Range("D2").Value = "String"
$wsCosts.Range("E2").Value = "Number"

# Rows 3/4 - field expression + search value, filled in column-by-column
# (Number column)
$wsCosts.Range("A3").Value = 'Field expression >>>'
$wsCosts.Range("A4").Value = 'Search values >>>'

$wsCosts.Range("B3").Value = 'Invoice ([a-z]{6}):;Number'
$wsCosts.Range("B4").Value = '(\d+-\d+);'

# (Date column)
$wsCosts.Range("C3").Value = 'Invoice ([a-z]{4}):;Date'
$wsCosts.Range("C4").Value = '([a-z]{3,}\s*\d+,\s*\d{4});'

# (Product column)
$wsCosts.Range("D3").Value = '[a-z]{7}:;Product'
$wsCosts.Range("D4").Value = '(.*oil.*);'

# (Total column)
$wsCosts.Range("E3").Value = '[a-z]{5}:;Total'
$wsCosts.Range("E4").Value = '\$([\d,]+\.*\d{2});'

# Column widths (best-fit-like, approximating the authored widths)
$wsCosts.Columns.Item(1).ColumnWidth = 18
$wsCosts.Columns.Item(2).ColumnWidth = 23.666666666666664
$wsCosts.Columns.Item(3).ColumnWidth = 11.166666666666668

# ---------------------------------------------------------------------------
# 4. Fix up "settlement" formatting / selection
# ---------------------------------------------------------------------------
# B3 / B4 pick up the same "Arial 10" style already used elsewhere on the sheet
$wsSettlement.Range("A3").Copy() | Out-Null
$wsSettlement.Range("B3:B4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsSettlement.Range("B1:D1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Final selection / active sheet -> "invoice_costs"
# ---------------------------------------------------------------------------
$wsCosts.Activate() | Out-Null
$wsCosts.Range("C3").Select() | Out-Null
